$d = $word.ActiveDocument

# Locate the run containing "https://www.nicolesharp.net/licenses/CC_BY-SA.htm"
# and collapse the search range down onto the matched text (no replacement
# performed by this call - it is only used to find the exact character
# offsets of the text we need to edit).
$rng = $d.Content
$found = $rng.Find.Execute("https://www.nicolesharp.net/licenses/CC_BY-SA.htm", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target URL text"
}

$start = $rng.Start
$end = $rng.End
$matched = $rng.Text

# The underscore sits right before "BY-SA.htm" at the end of the match.
$suffix = "BY-SA.htm"
$underscoreStart = $start + ($matched.Length - $suffix.Length - 1)
$underscoreEnd = $underscoreStart + 1

# Re-type the single "_" character as "-", in its own run, by toggling a
# character-formatting property around the text assignment so the run
# boundary is preserved instead of being silently re-merged with its
# neighbours on save.
$rHyphen = $d.Range($underscoreStart, $underscoreEnd)
$rHyphen.Bold = 1
$rHyphen.Text = "-"
$rHyphen.Bold = 0

# Re-type the trailing "BY-SA.htm" the same way so it lands in a run of
# its own, separate from both the hyphen run and the "...licenses/CC"
# prefix run.
$rSuffix = $d.Range($underscoreEnd, $end)
$rSuffix.Bold = 1
$rSuffix.Text = "~TMP~"
$rSuffix.Text = $suffix
$rSuffix.Bold = 0

$newEnd = $underscoreEnd + $suffix.Length
Write-Output "Updated URL text: $($d.Range($start, $newEnd).Text)"
